# Update crypto price/volume figures per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.713.09'
$ws.Range('E2').Value = '  -0.50%  '
$ws.Range('D3').Value = '2.462.20'
$ws.Range('E3').Value = '  +0.10%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').Value = '''573.14'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').Value = '''147.46'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  -1.64%  '
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('D11').Value = '''5.29'
$ws.Range('E11').Value = '  +0.18%  '
$ws.Range('D13').Value = '''29.07'
$ws.Range('E13').Value = '  +0.80%  '
$ws.Range('E14').Value = '  -1.67%  '
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '62.658.48'
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '2.465.75'
$ws.Range('E17').Value = '  +0.34%  '
$ws.Range('D18').Value = '''7.91'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').Value = '''10.90'
$ws.Range('E19').Value = '  -1.84%  '
$ws.Range('D20').Value = '''326.01'
$ws.Range('E20').Value = '  -1.31%  '
$ws.Range('E21').Value = '  -0.06%  '
$ws.Range('D22').Value = '''2.18'
$ws.Range('E22').Value = '  -2.68%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('D24').Value = '''9.97'
$ws.Range('E24').Value = '  +11.67%  '
$ws.Range('D25').Value = '''65.41'
$ws.Range('E25').Value = '  -1.59%  '
$ws.Range('D26').Value = '''640.56'
$ws.Range('E26').Value = '  -3.90%  '
$ws.Range('E27').Value = '  +0.21%  '
$ws.Range('D28').Value = '0.0₃0973'
$ws.Range('E28').Value = '  -2.83%  '
$ws.Range('E29').Value = '  -13.55%  '
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('E31').Value = '  -3.06%  '
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('E33').Value = '  -4.24%  '
$ws.Range('E34').Value = '  -0.08%  '
$ws.Range('E35').Value = '  -0.84%  '
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').Value = '''150.76'
$ws.Range('E38').Value = '  -1.53%  '
$ws.Range('D39').Value = '''18.56'
$ws.Range('E39').Value = '  -1.62%  '
$ws.Range('D40').Value = '''5.29'
$ws.Range('E40').Value = '  -3.92%  '
$ws.Range('D41').Value = '''2.71'
$ws.Range('E41').Value = '  -0.66%  '
$ws.Range('E42').Value = '  -2.27%  '
$ws.Range('E43').Value = '  -10.70%  '
$ws.Range('E44').Value = '  -0.05%  '
$ws.Range('D45').Value = '''153.06'
$ws.Range('E45').Value = '  +4.39%  '
$ws.Range('E46').Value = '  +1.12%  '
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('D48').Value = '''20.33'
$ws.Range('E48').Value = '  -1.98%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = '''0.0507'
$ws.Range('E50').Value = '  -1.74%  '
$ws.Range('E51').Value = '  -1.41%  '
